# Major update: re-sort the condition table (attendCond, A2:D13) ascending
# by column A ("attendCond"), the way Excel's Data > Sort dialog would,
# with "My data has headers" checked (header row A1:D1 stays put).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use the worksheet's Sort object (SortFields + Apply) rather than the bare
# Range.Sort shortcut so Excel records the sort as a reusable sortState
# (what the UI "Sort" dialog leaves behind in the saved sheet XML).
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A1:A13"), 0, 1, $null, 0)
$ws.Sort.SetRange($ws.Range("A1:D13"))
$ws.Sort.Header = 1
$ws.Sort.Apply()

# Leave the selection where the Sort dialog's "OK" click lands afterwards:
# the last key column's data cells.
[void]$ws.Range("A11:A13").Select()
